# FUNCTIONALITY: Relocated and wrote out a few test cases.
#
# The "CRUD" folder row (row 2) on Sheet1 pulls its Automated/Total test-case
# counts from a cached external-workbook link:
#   G2 = [1]Sheet1!$H$6   (Automated Test Cases)  202 -> 212
#   H2 = [1]Sheet1!$H$5   (Total Test Cases)       266 -> 276
#
# The source workbook (CRUD/_Test_Suite_Statistics_for_Folders.xlsx) isn't
# available in this environment, so we record the refreshed counts directly
# on the two cells that surface them. Every dependent formula (I2 = G2/H2,
# L5 = SUM($H:$H), L6 = SUM($G:$G), L7 = L6/L5) recalculates automatically
# from these new inputs, matching the rest of the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 276
$ws.Range("G2").Value = 212
